# Fruta / hortaliza, semanal
# Insert a new weekly record at row 105 (pushing the existing rows 105-110
# down to 106-111) on the only worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 105; this shifts rows
# 105..110 down to 106..111 and keeps everything else untouched.
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new data point.
$ws.Cells.Item(105, 1).Value  = 7
$ws.Cells.Item(105, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(105, 3).Value  = "Ñuble"
$ws.Cells.Item(105, 4).Value  = 44931
$ws.Cells.Item(105, 5).Value  = 16
$ws.Cells.Item(105, 6).Value  = 100112031
$ws.Cells.Item(105, 7).Value  = "Poroto verde"
$ws.Cells.Item(105, 8).Value  = "Sin especificar"
$ws.Cells.Item(105, 9).Value  = "Primera"
$ws.Cells.Item(105, 10).Value = 100
$ws.Cells.Item(105, 11).Value = 32000
$ws.Cells.Item(105, 12).Value = 33000
$ws.Cells.Item(105, 13).Value = 32500
$ws.Cells.Item(105, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(105, 15).Value = "Región del Maule"
$ws.Cells.Item(105, 16).Value = 1300
$ws.Cells.Item(105, 17).Value = 25
$ws.Cells.Item(105, 18).Value = "Hortaliza"
